$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.495
$ws.Range("B4").Value = 5.175

$ws.Range("A7").Value = -20.513

$ws.Range("B12").Value = 4.939

$ws.Range("A16").Value = -22.119

$ws.Range("B18").Value = 5.231
$ws.Range("B19").Value = 8.272
$ws.Range("B20").Value = 6.470999999999999

$ws.Range("A28").Value = -21.527
$ws.Range("A29").Value = -21.496

$ws.Range("B31").Value = 6.225999999999999

$ws.Range("A32").Value = -21.821

$ws.Range("A40").Value = -20.823
$ws.Range("B40").Value = 7.25

$ws.Range("B42").Value = 6.854000000000001

$ws.Range("B47").Value = 5.82
$ws.Range("B48").Value = 5.524

$ws.Range("A52").Value = -21.675

$ws.Range("A57").Value = -22.042

$ws.Range("B63").Value = 5.252000000000001
$ws.Range("B64").Value = 5.608000000000001

$ws.Range("A66").Value = -21.526

$ws.Range("B76").Value = 6.003

$ws.Range("B81").Value = 5.624

$ws.Range("B89").Value = 5.319999999999999

$ws.Range("B94").Value = 5.795000000000001

$ws.Range("A100").Value = -22.352
